$d = $word.ActiveDocument

# Word auto-creates the built-in "List Paragraph" style definition in
# styles.xml the first time it is referenced. Apply it momentarily to the
# existing content and then restore the original style so the visible
# document is unchanged, mirroring how Word mints this style on first use.
$p = $d.Paragraphs(1)
$originalStyle = $p.Range.Style
$p.Range.Style = "List Paragraph"
$p.Range.Style = $originalStyle

# Configure the newly minted style to match the standard built-in
# definition: based on Normal, ui priority 34, quick style, left indent
# 720 twips (36 pt) and contextual spacing (no space between paragraphs
# of the same style).
$listParagraph = $d.Styles("List Paragraph")
$listParagraph.BaseStyle = "Normal"
$listParagraph.Priority = 34
$listParagraph.QuickStyle = $true
$listParagraph.ParagraphFormat.LeftIndent = 36
$listParagraph.NoSpaceBetweenParagraphsOfSameStyle = $true
